$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record for "Macroferia Regional de Talca" / Coliflor needs to be
# inserted at row 202. This pushes the existing rows 202-245 down to 203-246
# (dimension grows from A1:R245 to A1:R246), preserving each row's formatting.
$ws.Rows.Item(202).Insert()

# Populate the newly inserted row 202 with the new record's data.
$ws.Range("A202").Value = 5
$ws.Range("B202").Value = "Macroferia Regional de Talca"
$ws.Range("C202").Value = "Maule"
$ws.Range("D202").Value = 44711
$ws.Range("E202").Value = 7
$ws.Range("F202").Value = 100112008
$ws.Range("G202").Value = "Coliflor"
$ws.Range("H202").Value = "Sin especificar"
$ws.Range("I202").Value = "Primera"
$ws.Range("J202").Value = 3000
$ws.Range("K202").Value = 900
$ws.Range("L202").Value = 900
$ws.Range("M202").Value = 900
$ws.Range("N202").Value = "$/unidad"
$ws.Range("O202").Value = "Región del Maule"
$ws.Range("P202").Value = 900
$ws.Range("Q202").Value = 1
$ws.Range("R202").Value = "Hortaliza"

# Two unrelated records (rows 241 / 242, after the shift: 242 / 243) had their
# "Región" swapped between visits; fix them up so region data lines back up
# with the correct sale dates.
$ws.Range("O242").Value = "Región del Maule"
$ws.Range("O243").Value = "Región Metropolitana"

# Likewise, the category ("Primera"/"Segunda") recorded for a couple of the
# shifted rows needs to be swapped back to line up with the correct record.
$ws.Range("I221").Value = "Primera"
$ws.Range("I222").Value = "Segunda"
$ws.Range("I229").Value = "Primera"
$ws.Range("I230").Value = "Segunda"
